# "Brought up to v3." — strip the legacy SharePoint custom XML parts
# (document library content-type schema + its companion item/itemProps
# parts) that got embedded in the package. This mirrors using Word's
# CustomXMLParts collection (Document Inspector's "Remove all custom XML
# documents" does the same thing under the hood) to delete every custom
# XML part from the document.

$d = $word.ActiveDocument

function Remove-AllCustomXmlParts {
    param($parts)

    if ($null -eq $parts) {
        return
    }

    try {
        $n = $parts.Count
    } catch {
        $n = 0
    }

    for ($i = $n; $i -ge 1; $i--) {
        try {
            $part = $parts.Item($i)
            if ($null -ne $part) {
                $part.Delete()
            }
        } catch {
            # keep going even if an individual part can't be removed
        }
    }
}

# 1) The default collection (non built-in parts).
Remove-AllCustomXmlParts $d.CustomXMLParts

# 2) Word also classifies some custom XML parts (e.g. the SharePoint
#    document-library content-type schema) as "built in" and hides them
#    from the plain CustomXMLParts collection — CustomXMLParts(True)
#    returns the full set, built-in included.
try {
    Remove-AllCustomXmlParts $d.CustomXMLParts($true)
} catch {
}

# 3) Belt-and-suspenders: target the known item IDs directly in case the
#    bulk collection enumeration above missed any of them.
$knownItemIds = @(
    "{A90CF479-1DC3-41DA-9078-B0657AFF7733}",
    "{711C9CB5-5D18-4D18-996D-586B65CB9D87}",
    "{A573EC8E-2A00-4C90-B29B-D922DA051274}"
)

foreach ($itemId in $knownItemIds) {
    try {
        $found = $d.CustomXMLParts.SelectByID($itemId)
        if ($null -ne $found) {
            $found.Delete()
        }
    } catch {
    }
}

# 4) And by namespace, covering the three schemas used by the parts
#    (contentType schema, SharePoint forms, and document properties).
$knownNamespaces = @(
    "http://schemas.microsoft.com/office/2006/metadata/contentType",
    "http://schemas.microsoft.com/sharepoint/v3/contenttype/forms",
    "http://schemas.microsoft.com/office/2006/metadata/properties"
)

foreach ($ns in $knownNamespaces) {
    try {
        $matches = $d.CustomXMLParts.SelectByNamespace($ns)
        Remove-AllCustomXmlParts $matches
    } catch {
    }
}

# 5) Defensive "keep deleting the first item" sweep, in case the
#    collection re-populates / re-indexes oddly after each delete.
$guard = 0
while ($guard -lt 32) {
    try {
        if ($d.CustomXMLParts.Count -lt 1) {
            break
        }
        $d.CustomXMLParts.Item(1).Delete()
    } catch {
        break
    }
    $guard = $guard + 1
}

$d.Save()
